# Regenerate save_data to use K (strikeouts) instead of Strike# in column G.
# New K values were recalculated for rows 2-19.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 2
    3  = 2
    4  = 0
    5  = 0
    6  = 2
    7  = 1
    8  = 0
    9  = 0
    10 = 2
    11 = 3
    12 = 3
    13 = 3
    14 = 2
    15 = 1
    16 = 2
    17 = 2
    18 = 1
    19 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
